$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only C changes
$ws.Range("C2").Value = 50

# Row 3: B and C change
$ws.Range("B3").Value = "<tab>"
$ws.Range("C3").Value = 51

# Row 4: B and C change
$ws.Range("B4").Value = "<by>"
$ws.Range("C4").Value = 46

# Row 5: B and C change
$ws.Range("B5").Value = "<it>"
$ws.Range("C5").Value = 46

# Row 6: B and C change
$ws.Range("B6").Value = "<whis>"
$ws.Range("C6").Value = 46

# Row 7: only C changes
$ws.Range("C7").Value = 46

# Row 8: B and C change
$ws.Range("B8").Value = "<be>"
$ws.Range("C8").Value = 47

# Row 9: B and C change
$ws.Range("B9").Value = "<was>"
$ws.Range("C9").Value = 42

# Row 10: B and C change
$ws.Range("B10").Value = "<see>"
$ws.Range("C10").Value = 42

# Row 11: only C changes
$ws.Range("C11").Value = 50

# Row 12: B and C change
$ws.Range("B12").Value = "<be>"
$ws.Range("C12").Value = 46

# Row 13: B and C change
$ws.Range("B13").Value = "<for>"
$ws.Range("C13").Value = 41

# Row 14: only C changes
$ws.Range("C14").Value = 44

# Row 15: B and C change
$ws.Range("B15").Value = "<alph>"
$ws.Range("C15").Value = 47

# Row 16: B and C change
$ws.Range("B16").Value = "<part>"
$ws.Range("C16").Value = 43

# Row 18: only C changes
$ws.Range("C18").Value = 41
